$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new (blank) column before column N, shifting the existing
# Late/Date/Outstanding columns (N:P) one column to the right (O:Q).
$ws.Columns("N").Insert()

# Give the freshly inserted column the same on-screen width as column M
# (11 characters), matching the workbook author's manual resize.
$ws.Columns("N").ColumnWidth = 11 - 5/6

# Activate the "Repayment schedule" sheet and move the selection there,
# matching the recorded cursor position after the edit.
$ws.Activate() | Out-Null
$ws.Range("I21").Select() | Out-Null
